$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01473133333333333
$ws.Range("H2").Value = 0.044194
$ws.Range("I2").Value = [double]"4.600893861834377E-05"
$ws.Range("J2").Value = [double]"4.600893861834377E-05"
$ws.Range("M2").Value = 13.35941066666667
$ws.Range("N2").Value = 40.078232
$ws.Range("O2").Value = 0.4925555025958562
$ws.Range("P2").Value = 0.4925555025958562
$ws.Range("Q2").Value = 0.1968019316675556
$ws.Range("R2").Value = 1.771217385008
$ws.Range("S2").Value = [double]"2.266195588506021E-05"
$ws.Range("T2").Value = [double]"2.266195588506021E-05"

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01473133333333333
$ws.Range("H3").Value = 0.044194
$ws.Range("I3").Value = [double]"4.600893861834377E-05"
$ws.Range("J3").Value = [double]"4.600893861834377E-05"
$ws.Range("O3").Value = 0.03774352140193379
$ws.Range("P3").Value = 0.03774352140193379
$ws.Range("Q3").Value = 0.01508052976911111
$ws.Range("R3").Value = 0.135724767922
$ws.Range("S3").Value = [double]"1.736539359421716E-06"
$ws.Range("T3").Value = [double]"1.736539359421716E-06"

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01473133333333333
$ws.Range("H4").Value = 0.044194
$ws.Range("I4").Value = [double]"4.600893861834377E-05"
$ws.Range("J4").Value = [double]"4.600893861834377E-05"
$ws.Range("M4").Value = 12.73953533333333
$ws.Range("N4").Value = 38.218606
$ws.Range("O4").Value = 0.4697009760022101
$ws.Range("P4").Value = 0.46970097600221
$ws.Range("Q4").Value = 0.1876703415071111
$ws.Range("R4").Value = 1.689033073564
$ws.Range("S4").Value = [double]"2.161044337386184E-05"
$ws.Range("T4").Value = [double]"2.161044337386184E-05"

# Row 5
$ws.Range("I5").Value = 0.9987490355524334
$ws.Range("J5").Value = 0.9987490355524334
$ws.Range("M5").Value = 13.35941066666667
$ws.Range("N5").Value = 40.078232
$ws.Range("O5").Value = 0.4925555025958562
$ws.Range("P5").Value = 0.4925555025958562
$ws.Range("Q5").Value = 4272.120708506417
$ws.Range("R5").Value = 38449.08637655775
$ws.Range("S5").Value = 0.4919393331736555
$ws.Range("T5").Value = 0.4919393331736555

# Row 6
$ws.Range("I6").Value = 0.9987490355524334
$ws.Range("J6").Value = 0.9987490355524334
$ws.Range("O6").Value = 0.03774352140193379
$ws.Range("P6").Value = 0.03774352140193379
$ws.Range("S6").Value = 0.03769630559853401
$ws.Range("T6").Value = 0.037696305598534

# Row 7
$ws.Range("I7").Value = 0.9987490355524334
$ws.Range("J7").Value = 0.9987490355524334
$ws.Range("M7").Value = 12.73953533333333
$ws.Range("N7").Value = 38.218606
$ws.Range("O7").Value = 0.4697009760022101
$ws.Range("P7").Value = 0.46970097600221
$ws.Range("Q7").Value = 4073.894730257752
$ws.Range("R7").Value = 36665.05257231976
$ws.Range("S7").Value = 0.469113396780244
$ws.Range("T7").Value = 0.4691133967802439

# Row 8
$ws.Range("G8").Value = 0.3858076666666667
$ws.Range("H8").Value = 1.157423
$ws.Range("I8").Value = 0.001204955508948258
$ws.Range("J8").Value = 0.001204955508948258
$ws.Range("M8").Value = 13.35941066666667
$ws.Range("N8").Value = 40.078232
$ws.Range("O8").Value = 0.4925555025958562
$ws.Range("P8").Value = 0.4925555025958562
$ws.Range("Q8").Value = 5.154163057348446
$ws.Range("R8").Value = 46.38746751613601
$ws.Range("S8").Value = 0.0005935074663156549
$ws.Range("T8").Value = 0.0005935074663156548

# Row 9
$ws.Range("G9").Value = 0.3858076666666667
$ws.Range("H9").Value = 1.157423
$ws.Range("I9").Value = 0.001204955508948258
$ws.Range("J9").Value = 0.001204955508948258
$ws.Range("O9").Value = 0.03774352140193379
$ws.Range("P9").Value = 0.03774352140193379
$ws.Range("Q9").Value = 0.3949529801998889
$ws.Range("R9").Value = 3.554576821799
$ws.Range("S9").Value = [double]"4.54792640403666E-05"
$ws.Range("T9").Value = [double]"4.547926404036659E-05"

# Row 10
$ws.Range("G10").Value = 0.3858076666666667
$ws.Range("H10").Value = 1.157423
$ws.Range("I10").Value = 0.001204955508948258
$ws.Range("J10").Value = 0.001204955508948258
$ws.Range("M10").Value = 12.73953533333333
$ws.Range("N10").Value = 38.218606
$ws.Range("O10").Value = 0.4697009760022101
$ws.Range("P10").Value = 0.46970097600221
$ws.Range("Q10").Value = 4.91501040137089
$ws.Range("R10").Value = 44.23509361233801
$ws.Range("S10").Value = 0.0005659687785922365
$ws.Range("T10").Value = 0.0005659687785922364
